$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Total Points" row - Nick's total changed from 150 to 163
$ws.Range("B2").Value = 163

# "Rolls" row - new roll history strings for each player
$ws.Range("B6").Value = "0/0/0/0/"
$ws.Range("C6").Value = "1/1/"
$ws.Range("D6").Value = "2/"
$ws.Range("E6").Value = "3/3/3/3/3/"
